# Update cryptos list values (price + 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "68.246.29"
Set-TextValue "E2" "  -0.36%  "
Set-TextValue "D3" "2.647.18"
Set-TextValue "E3" "  +0.07%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "597.49"
Set-TextValue "E5" "  -0.43%  "
Set-TextValue "D6" "156.64"
Set-TextValue "E6" "  +1.15%  "
Set-TextValue "E8" "  -0.38%  "
Set-TextValue "E9" "  +2.73%  "
Set-TextValue "E10" "  -1.20%  "
Set-TextValue "E11" "  +0.56%  "
Set-TextValue "E12" "  +0.63%  "
Set-TextValue "D13" "28.01"
Set-TextValue "E13" "  +0.58%  "
Set-TextValue "E14" "  +1.31%  "
Set-TextValue "D15" "3.129.35"
Set-TextValue "E15" "  +0.13%  "
Set-TextValue "D16" "68.346.90"
Set-TextValue "E16" "  -0.06%  "
Set-TextValue "D17" "2.649.26"
Set-TextValue "E17" "  +0.13%  "
Set-TextValue "E18" "  -0.37%  "
Set-TextValue "D19" "363.44"
Set-TextValue "E19" "  -1.30%  "
Set-TextValue "E20" "  -1.17%  "
Set-TextValue "E21" "  +3.26%  "
Set-TextValue "E22" "  -1.14%  "
Set-TextValue "E23" "  -2.44%  "
Set-TextValue "E24" "  +2.59%  "
Set-TextValue "E25" "  +0.01%  "
Set-TextValue "D26" "9.76"
Set-TextValue "E26" "  -2.53%  "
Set-TextValue "D27" "2.779.55"
Set-TextValue "E27" "  +0.32%  "
Set-TextValue "E28" "  -0.77%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.08%  "
Set-TextValue "D30" "559.99"
Set-TextValue "E30" "  -2.76%  "
Set-TextValue "E31" "  +0.66%  "
Set-TextValue "E32" "  -0.94%  "
Set-TextValue "E33" "  +0.35%  "
Set-TextValue "E34" "  -0.43%  "
Set-TextValue "E35" "  +0.02%  "
Set-TextValue "E36" "  +0.50%  "
Set-TextValue "D37" "161.89"
Set-TextValue "E37" "  +1.95%  "
Set-TextValue "D38" "19.68"
Set-TextValue "E38" "  +2.38%  "
Set-TextValue "D39" "0.371"
Set-TextValue "E39" "  +1.16%  "
Set-TextValue "E40" "  -2.84%  "
Set-TextValue "E41" "  -1.03%  "
Set-TextValue "D42" "0.0₆0335"
Set-TextValue "E42" "  +3.85%  "
Set-TextValue "D43" "17.79"
Set-TextValue "E43" "  +0.26%  "
Set-TextValue "E44" "  -1.66%  "
Set-TextValue "E45" "  +0.03%  "
Set-TextValue "D46" "158.89"
Set-TextValue "D47" "3.72"
Set-TextValue "E47" "  -0.47%  "
Set-TextValue "D48" "22.07"
Set-TextValue "E48" "  +0.55%  "
Set-TextValue "E49" "  -1.35%  "
Set-TextValue "E50" "  +0.42%  "
Set-TextValue "D51" "0.614"
Set-TextValue "E51" "  -0.27%  "
